$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text; some look like numbers (e.g. "162.00")
# so force text format first to avoid Excel auto-converting them to numeric values
# and dropping significant trailing zeros / precision.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.985.53"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.398.86"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.14"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.00"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.396.33"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.420"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.991.67"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.78"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.025.09"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.416.21"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.38"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.13"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.76"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.28"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.42"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.04"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.38"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.73"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.00"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.73"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.853"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.81"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0722"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.70"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.61"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.729.45"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.39"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.77"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.33"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0302"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "326.80"

# Column E (Volume/1h change) values keep their leading/trailing double-space padding
# so Excel keeps them as text automatically.
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  -5.49%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  -4.63%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("E25").Value = "  -3.91%  "
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("E27").Value = "  -4.98%  "
$ws.Range("E28").Value = "  -1.16%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("E36").Value = "  -7.54%  "
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  +6.62%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -5.26%  "
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  -2.27%  "
